$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.08390303802875526
$ws.Range("H2").Value = -36.15692857716881
$ws.Range("I2").Value = -21.20859238905999
$ws.Range("G3").Value = 0.1299686832758773
$ws.Range("H3").Value = 45.97122178476677
$ws.Range("G4").Value = -0.763411951047763
$ws.Range("H4").Value = -24.06274390632854
$ws.Range("G5").Value = -0.6487145077394124
$ws.Range("H5").Value = -6.243315545876628
$ws.Range("G6").Value = 0.2327571647385957
$ws.Range("H6").Value = -5.424704002520047
$ws.Range("G7").Value = 0.3233338588108007
$ws.Range("H7").Value = 97.36907851642994
$ws.Range("G8").Value = 0.1438878313958999
$ws.Range("H8").Value = -12.95822143315205
$ws.Range("G9").Value = 0.1833736205203756
$ws.Range("H9").Value = -6.006851129854972
$ws.Range("G10").Value = -0.1561912374687073
$ws.Range("H10").Value = -173.3179703324856
$ws.Range("G11").Value = -0.125363707273024
$ws.Range("H11").Value = -5.553150946702916
$ws.Range("G12").Value = 0.2174705154782628
$ws.Range("H12").Value = 36.74872171655078
$ws.Range("G13").Value = 0.1798461969837346
$ws.Range("H13").Value = -12.55243167736414
$ws.Range("G14").Value = 0.190004964532156
$ws.Range("H14").Value = 0.3336737373904686
$ws.Range("G15").Value = 0.2367767550015582
$ws.Range("H15").Value = -5.253452582987111
$ws.Range("G16").Value = 0.005022356227122955
$ws.Range("H16").Value = -86.23266930830148
$ws.Range("G17").Value = -0.0001728912089638381
$ws.Range("H17").Value = -100.4874212938964
$ws.Range("G18").Value = 0.08548085485803285
$ws.Range("H18").Value = -50.67868453485465
$ws.Range("G19").Value = 0.104218744917061
$ws.Range("H19").Value = -17.12421059670918
$ws.Range("G20").Value = 0.08767697570785186
$ws.Range("H20").Value = -23.52774038717628
$ws.Range("G21").Value = 0.08351784939982017
$ws.Range("H21").Value = -16.81391827364226
$ws.Range("G22").Value = 0.07631341488650217
$ws.Range("H22").Value = -18.98541327627335
$ws.Range("G23").Value = 0.05666116387261208
$ws.Range("H23").Value = -47.7729097528159
$ws.Range("G24").Value = -0.2338513118146741
$ws.Range("H24").Value = -87.64496401943317
$ws.Range("G25").Value = -0.1826443741354775
$ws.Range("H25").Value = 17.89569846713929
$ws.Range("G26").Value = 0.1837546816765711
$ws.Range("H26").Value = 15.58219947580966
$ws.Range("G27").Value = 0.2045339896368011
$ws.Range("H27").Value = 2.029095406105847
$ws.Range("G28").Value = -0.02339979522631088
$ws.Range("H28").Value = -191.3883358039575
$ws.Range("G29").Value = 0.03358686589417845
$ws.Range("H29").Value = 118.4157613836284
